# Add 2022-Q4 data
# 1) Insert a new worksheet named "2022-Q4" right after "2022-Q3" position,
#    i.e. right before the current second sheet (so it becomes sheet index 2,
#    pushing 2022-Q3/2021-Q4/2021-Q3/2020-Q4 down by one).
# 2) Populate the "总计" (summary) sheet with a new row for 2022-Q4 at the
#    top of the data (row 2), shifting the other quarters down by one row.
# 3) Populate the new "2022-Q4" sheet with the per-fund detail data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: Update the "总计" summary sheet (sheet index 1)
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Shift existing rows 2-5 down to rows 3-6, copying formatting along the way
# (copy from bottom to top so we don't overwrite rows we still need to read).
# Only copy the B:D columns (the text/number data); the A column is just a
# sequential 0-based row index so it is rewritten explicitly afterwards.
$summary.Range("B5:D5").Copy($summary.Range("B6:D6"))
$summary.Range("B4:D4").Copy($summary.Range("B5:D5"))
$summary.Range("B3:D3").Copy($summary.Range("B4:D4"))
$summary.Range("B2:D2").Copy($summary.Range("B3:D3"))
$summary.Range("A5").Copy($summary.Range("A6"))

# Fill in the new row 2 with the 2022-Q4 figures
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 1.49

# Rewrite the A column's sequential 0-based index for every row
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------
# Step 2: Insert the new "2022-Q4" worksheet before the current "2022-Q3"
#          worksheet (sheet index 2). We clone the existing "2022-Q3"
#          sheet (same sheetPr/column layout/header/styles) and place the
#          clone immediately before it, then overwrite the clone's data
#          with the 2022-Q4 figures.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# "2022-Q3" has 5 fund rows (rows 2-6); "2022-Q4" only has 4 fund rows
# (rows 2-5), so remove the now-extraneous last row that was copied over.
$q4.Rows.Item(6).Delete()

# The header (row 1) text is identical across quarters, so it is left
# untouched; only the fund data rows (2-5) are overwritten below.

# Numeric-looking text (fund codes with leading zeros, decimal figures
# stored as text) must be forced to the "Text" number format before the
# value is assigned - otherwise Excel auto-converts them to real numbers
# (dropping leading zeros, etc). Plain (non-numeric-looking) text, like the
# fund name, does not need this treatment.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2
Set-TextValue $q4.Range("B2") "010695"
$q4.Range("C2").Value = "华夏磐益一年定期开放混合"
Set-TextValue $q4.Range("D2") "16.03"
Set-TextValue $q4.Range("E2") "98.69"
Set-TextValue $q4.Range("F2") "4.53"
Set-TextValue $q4.Range("G2") "0.7262"
$q4.Range("H2").Value = 6

# Row 3
Set-TextValue $q4.Range("B3") "004702"
$q4.Range("C3").Value = "南方金融主题灵活配置混合A"
Set-TextValue $q4.Range("D3") "12.97"
Set-TextValue $q4.Range("E3") "92.71"
Set-TextValue $q4.Range("F3") "4.19"
Set-TextValue $q4.Range("G3") "0.5434"
$q4.Range("H3").Value = 7

# Row 4
Set-TextValue $q4.Range("B4") "013500"
$q4.Range("C4").Value = "南方金融主题灵活配置混合C"
Set-TextValue $q4.Range("D4") "4.80"
Set-TextValue $q4.Range("E4") "92.71"
Set-TextValue $q4.Range("F4") "4.19"
Set-TextValue $q4.Range("G4") "0.2011"
$q4.Range("H4").Value = 7

# Row 5
Set-TextValue $q4.Range("B5") "000270"
$q4.Range("C5").Value = "建信灵活配置混合"
Set-TextValue $q4.Range("D5") "1.53"
Set-TextValue $q4.Range("E5") "91.22"
Set-TextValue $q4.Range("F5") "1.01"
Set-TextValue $q4.Range("G5") "0.0155"
$q4.Range("H5").Value = 3

Write-Output "done"
